$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "60.002.41"
Set-TextValue $ws.Range("E2") "  +2.22%  "
Set-TextValue $ws.Range("D3") "3.187.87"
Set-TextValue $ws.Range("E3") "  +0.86%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "536.55"
Set-TextValue $ws.Range("E5") "  +1.26%  "
Set-TextValue $ws.Range("D6") "144.97"
Set-TextValue $ws.Range("E6") "  +3.64%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("E8") "  -2.44%  "
Set-TextValue $ws.Range("E9") "  +0.33%  "
Set-TextValue $ws.Range("E10") "  +0.63%  "
Set-TextValue $ws.Range("E11") "  -1.82%  "
Set-TextValue $ws.Range("D12") "3.737.47"
Set-TextValue $ws.Range("E12") "  +0.88%  "
Set-TextValue $ws.Range("D13") "0.137"
Set-TextValue $ws.Range("E13") "  -2.82%  "
Set-TextValue $ws.Range("D14") "25.74"
Set-TextValue $ws.Range("E14") "  -0.76%  "
Set-TextValue $ws.Range("E15") "  -0.34%  "
Set-TextValue $ws.Range("D16") "60.018.07"
Set-TextValue $ws.Range("E16") "  +2.18%  "
Set-TextValue $ws.Range("D17") "3.183.71"
Set-TextValue $ws.Range("E17") "  +1.50%  "
Set-TextValue $ws.Range("E18") "  -0.60%  "
Set-TextValue $ws.Range("E19") "  +1.62%  "
Set-TextValue $ws.Range("E20") "  +0.58%  "
Set-TextValue $ws.Range("D21") "368.64"
Set-TextValue $ws.Range("E21") "  -1.82%  "
Set-TextValue $ws.Range("E22") "  +0.02%  "
Set-TextValue $ws.Range("E23") "  -2.08%  "
Set-TextValue $ws.Range("E24") "  -0.30%  "
Set-TextValue $ws.Range("D25") "0.170"
Set-TextValue $ws.Range("E25") "  +1.70%  "
Set-TextValue $ws.Range("E26") "  +3.52%  "
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  +0.08%  "
Set-TextValue $ws.Range("D28") "0.0₃0873"
Set-TextValue $ws.Range("E28") "  +0.50%  "
Set-TextValue $ws.Range("D29") "22.43"
Set-TextValue $ws.Range("E29") "  +0.37%  "
Set-TextValue $ws.Range("E30") "  +0.31%  "
Set-TextValue $ws.Range("D31") "6.10"
Set-TextValue $ws.Range("E31") "  +0.39%  "
Set-TextValue $ws.Range("D32") "5.27"
Set-TextValue $ws.Range("E32") "  +2.71%  "
Set-TextValue $ws.Range("E33") "  +4.89%  "
Set-TextValue $ws.Range("E34") "  +2.43%  "
Set-TextValue $ws.Range("D35") "157.74"
Set-TextValue $ws.Range("E35") "  -0.37%  "
Set-TextValue $ws.Range("D36") "1.36"
Set-TextValue $ws.Range("E36") "  +1.53%  "
Set-TextValue $ws.Range("D37") "26.26"
Set-TextValue $ws.Range("E37") "  +5.10%  "
Set-TextValue $ws.Range("D38") "2.784.38"
Set-TextValue $ws.Range("E38") "  +5.13%  "
Set-TextValue $ws.Range("D39") "0.0709"
Set-TextValue $ws.Range("E39") "  +2.86%  "
Set-TextValue $ws.Range("B40") "VeChain"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.0307"
Set-TextValue $ws.Range("E40") "  +6.17%  "
Set-TextValue $ws.Range("B41") "Stacks"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "1.69"
Set-TextValue $ws.Range("E41") "  -0.06%  "
Set-TextValue $ws.Range("E42") "  -1.83%  "
Set-TextValue $ws.Range("D43") "39.87"
Set-TextValue $ws.Range("E43") "  +1.92%  "
Set-TextValue $ws.Range("E44") "  -0.38%  "
Set-TextValue $ws.Range("E45") "  +0.92%  "
Set-TextValue $ws.Range("D46") "3.228.81"
Set-TextValue $ws.Range("E46") "  +0.88%  "
Set-TextValue $ws.Range("D47") "0.981"
Set-TextValue $ws.Range("E47") "  +0.34%  "
Set-TextValue $ws.Range("E48") "  -0.95%  "
Set-TextValue $ws.Range("D49") "0.795"
Set-TextValue $ws.Range("E49") "  +5.77%  "
Set-TextValue $ws.Range("D50") "20.55"
Set-TextValue $ws.Range("E50") "  +2.33%  "
Set-TextValue $ws.Range("E51") "  +0.02%  "

Write-Output "Updated cryptos list"
